$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 new values
$ws.Range("A23").Value = 112026905
$ws.Range("B23").Value = 90799
$ws.Range("E23").Value = 1968
$ws.Range("F23").Value = "Grantaggsvamp"
$ws.Range("G23").Value = "Bankera violascens"
$ws.Range("H23").Value = "(Alb. & Schwein. : Fr.) Pouzar"
$ws.Range("Q23").Value = 485427
$ws.Range("R23").Value = 6996682

# Row 24 new values
$ws.Range("A24").Value = 112026957
$ws.Range("B24").Value = 85448
$ws.Range("E24").Value = 3739
$ws.Range("F24").Value = "Persiljespindling"
$ws.Range("G24").Value = "Cortinarius sulfurinus"
$ws.Range("H24").Value = "Quél."
$ws.Range("Q24").Value = 485421
$ws.Range("R24").Value = 6996666
